# Apply cryptocurrency price/volume updates per the Thu Jan 11 17:32:48 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.816.69'
$ws.Range('E2').Value = '  +0.42%  '

$ws.Range('D3').Value = '2.580.31'
$ws.Range('E3').Value = '  +6.77%  '

$ws.Range('E4').Value = '  +0.22%  '

$ws.Range('D5').Value = '''305.83'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.92%  '

$ws.Range('D6').Value = '''98.78'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.43%  '

$ws.Range('D7').Value = '''0.598'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +5.33%  '

$ws.Range('E8').Value = '  +0.25%  '

$ws.Range('D9').Value = '''0.574'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +14.30%  '

$ws.Range('D10').Value = '''38.85'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +11.77%  '

$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').Value = '''54.09'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.95%  '

$ws.Range('B12').Value = 'Dogecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D12').Value = '''0.0839'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.72%  '

$ws.Range('D13').Value = '''8.12'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +15.31%  '

$ws.Range('D14').Value = '2.995.11'
$ws.Range('E14').Value = '  +7.29%  '

$ws.Range('D15').Value = '''0.105'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.10%  '

$ws.Range('D16').Value = '2.601.86'
$ws.Range('E16').Value = '  +7.80%  '

$ws.Range('D17').Value = '''0.917'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +9.35%  '

$ws.Range('D18').Value = '''14.90'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +5.52%  '

$ws.Range('D19').Value = '46.131.92'
$ws.Range('E19').Value = '  +1.27%  '

$ws.Range('D20').Value = '''0.0000100'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +6.53%  '

$ws.Range('D21').Value = '''12.89'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.97%  '

$ws.Range('D22').Value = '''6.68'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +8.45%  '

$ws.Range('D23').Value = '''71.07'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.94%  '

$ws.Range('D24').Value = '''270.65'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +11.98%  '

$ws.Range('D25').Value = '''3.01'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +8.22%  '

$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '''2.15'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +11.77%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''29.65'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +39.89%  '

$ws.Range('D28').Value = '''1.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.03%  '

$ws.Range('D29').Value = '''4.02'
$ws.Range('D29').Style = "Normal"

$ws.Range('D30').Value = '''10.52'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +8.82%  '

$ws.Range('D31').Value = '''2.31'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.48%  '

$ws.Range('D32').Value = '''38.68'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.30%  '

$ws.Range('E33').Value = '  +14.57%  '

$ws.Range('D34').Value = '''3.61'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.92%  '

$ws.Range('E35').Value = '  +2.46%  '

$ws.Range('D36').Value = '''0.0834'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +8.92%  '

$ws.Range('D37').Value = '''2.18'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +9.68%  '

$ws.Range('D38').Value = '''149.24'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.14%  '

$ws.Range('D39').Value = '''0.120'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +6.12%  '

$ws.Range('D40').Value = '''0.121'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.01%  '

$ws.Range('D41').Value = '''23.04'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +43.10%  '

$ws.Range('D42').Value = '''15.81'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +8.55%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0327'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +9.90%  '

$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = '''3.58'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +11.72%  '

$ws.Range('D45').Value = '''4.06'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +7.84%  '

$ws.Range('D46').Value = '2.152.22'
$ws.Range('E46').Value = '  +8.23%  '

$ws.Range('E47').Value = '  +0.07%  '

$ws.Range('D48').Value = '''92.93'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.59%  '

$ws.Range('D49').Value = '''9.54'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +12.04%  '

$ws.Range('D50').Value = '''108.45'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +7.85%  '

$ws.Range('D51').Value = '''1.76'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.87%  '
